$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / workbook view tweak -------------------------------------------
# Widen the sheet-tab area (tabRatio) left by the tab scroll, as recorded by
# Excel after editing the sheet.
$excel.ActiveWindow.TabRatio = 0.854

# --- New data-driven test rows (DataDriver library) -----------------------
# Row 2: add a second "username" test value (Bob) next to existing "Admin"
$ws.Range("B2").Value = "Bob"

# Row 4: add a second "password" test value (Jane) next to existing "admin123"
$ws.Range("A4").Value = "Jane"

# --- Column widths ----------------------------------------------------------
# Widen columns A:B so the longer test data/credentials are readable
$ws.Columns("A:B").ColumnWidth = 30.71

# --- Marker cells used by DataDriver (empty but carry formatting) ----------
# A5: alignment-flavoured style (reading order) placeholder row after data
$ws.Range("A5").HorizontalAlignment = -4108

# R8: protection-flavoured style placeholder cell (unlocked) far right
$ws.Range("R8").Locked = $false

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.Orientation = 1
